$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 229. This pushes the existing
# rows 229:246 down to 231:246, preserving all their values/formatting,
# and grows the sheet's used range to A1:R248 (matches the new <dimension>).
$ws.Rows("229:230").Insert()

# Populate the newly inserted row 229 with the new weekly record.
$ws.Range("A229").Value = 7
$ws.Range("B229").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C229").Value = "Ñuble"
$ws.Range("D229").Value = 44769
$ws.Range("E229").Value = 16
$ws.Range("F229").Value = 100112006
$ws.Range("G229").Value = "Repollo"
$ws.Range("H229").Value = "Crespo record"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 200
$ws.Range("K229").Value = 1200
$ws.Range("L229").Value = 1400
$ws.Range("M229").Value = 1300
$ws.Range("N229").Value = "`$/unidad"
$ws.Range("O229").Value = "Provincia de Diguillín"
$ws.Range("P229").Value = 1300
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"

# Populate the newly inserted row 230 with the new weekly record.
$ws.Range("A230").Value = 7
$ws.Range("B230").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C230").Value = "Ñuble"
$ws.Range("D230").Value = 44769
$ws.Range("E230").Value = 16
$ws.Range("F230").Value = 100112006
$ws.Range("G230").Value = "Repollo"
$ws.Range("H230").Value = "Crespo record"
$ws.Range("I230").Value = "Segunda"
$ws.Range("J230").Value = 150
$ws.Range("K230").Value = 1000
$ws.Range("L230").Value = 1000
$ws.Range("M230").Value = 1000
$ws.Range("N230").Value = "`$/unidad"
$ws.Range("O230").Value = "Provincia de Diguillín"
$ws.Range("P230").Value = 1000
$ws.Range("Q230").Value = 1
$ws.Range("R230").Value = "Hortaliza"
